$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the old email address with the placeholder email address
$ws.Range("A2").Value = "EnterYourEmailAddress@Here.com"
$ws.Range("A3").Value = "EnterYourEmailAddress@Here.com"

$ws.Range("A3").Select()
